# Add Bruce Kendall's name to column A (CompadrinoName) for data rows 12-43.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 12; $r -le 43; $r++) {
    $ws.Cells.Item($r, 1).Value = "Bruce Kendall"
}

# Restore the selection/view position to reflect where the author ended up.
$ws.Range("C41").Select()
